$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: updated "datos actualizados" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 15:04"

# Row 4: Estados Unidos - value updates only
$ws.Range("B4").Value = 2891380
$ws.Range("C4").Value = 792
$ws.Range("E4").Value = 1523303

# Row 7: India - value updates only
$ws.Range("B7").Value = 651065
$ws.Range("C7").Value = 1176
$ws.Range("D7").Value = 395128
$ws.Range("E7").Value = 237242
$ws.Range("G7").Value = 26
$ws.Range("H7").Value = 18695

# Row 14: Iran - value updates only (keeps its rank)
$ws.Range("B14").Value = 237878
$ws.Range("C14").Value = 2449
$ws.Range("D14").Value = 198949
$ws.Range("E14").Value = 27521
$ws.Range("G14").Value = 148
$ws.Range("H14").Value = 11408

# Rows 16-17: Arabia Saudita overtakes Turquia in total cases, so they swap places
$ws.Range("A16").Value = "Arabia Saudita"
$ws.Range("B16").Value = 205929
$ws.Range("C16").Value = 4128
$ws.Range("D16").Value = 143256
$ws.Range("E16").Value = 60815
$ws.Range("G16").Value = 56
$ws.Range("H16").Value = 1858

$ws.Range("A17").Value = "Turquia"
$ws.Range("B17").Value = 203456
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 178278
$ws.Range("E17").Value = 19992
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 5186

# Row 48: Suiza - value updates only
$ws.Range("B48").Value = 32198
$ws.Range("C48").Value = 97
$ws.Range("E48").Value = 1033

# Rows 101-102: Croacia overtakes Guinea Ecuatorial in total cases, so they swap places
$ws.Range("A101").Value = "Croacia"
$ws.Range("B101").Value = 3094
$ws.Range("C101").Value = 86
$ws.Range("D101").Value = 2183
$ws.Range("E101").Value = 798
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 113

$ws.Range("A102").Value = "Guinea Ecuatorial"
$ws.Range("B102").Value = 3071
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 842
$ws.Range("E102").Value = 2178
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 51

# Row 104: Albania - value updates only
$ws.Range("B104").Value = 2819
$ws.Range("C104").Value = 67
$ws.Range("D104").Value = 1637
$ws.Range("E104").Value = 1108
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 74

# Rows 105-107: Madagascar overtakes Mayotte and Nicaragua in total cases, so it moves up two spots
$ws.Range("A105").Value = "Madagascar"
$ws.Range("B105").Value = 2728
$ws.Range("C105").Value = 216
$ws.Range("D105").Value = 1078
$ws.Range("E105").Value = 1621
$ws.Range("G105").Value = 3
$ws.Range("H105").Value = 29

$ws.Range("A106").Value = "Mayotte"
$ws.Range("B106").Value = 2661
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 2375
$ws.Range("E106").Value = 251
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 35

$ws.Range("A107").Value = "Nicaragua"
$ws.Range("B107").Value = 2519
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 1238
$ws.Range("E107").Value = 1198
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 83
